{"js": "// This script updates the worksheet/body date line and the 100 arithmetic\n// practice cells in the table. Every one of the 101 paragraphs in the\n// document body (the heading date paragraph followed by the 100 table\n// cell paragraphs, in document order) gets its text replaced with the\n// corresponding new value below. We operate strictly by paragraph index\n// (document order) rather than by searching for old text, since some of\n// the new values collide with other cells' old values (e.g. \"72-5=\" is\n// an old value for one cell and the new value for another), which would\n// make a naive find/replace approach ambiguous or order-dependent.\n\nconst newValues = [\"2025-01-02 Thursday\", \"5+12=\", \"58-0=\", \"62+20=\", \"78-75=\", \"88+5=\", \"17+74=\", \"97-29=\", \"65-28=\", \"27+44=\", \"39+57=\", \"50-5=\", \"28+39=\", \"16+25=\", \"3+91=\", \"31-31=\", \"37+45=\", \"10+32=\", \"1+23=\", \"65-49=\", \"60-48=\", \"18+17=\", \"42+17=\", \"42+17=\", \"91+3=\", \"27+33=\", \"11+25=\", \"92-50=\", \"52+17=\", \"43+15=\", \"9+24=\", \"52-10=\", \"48-22=\", \"62+15=\", \"5+35=\", \"16+70=\", \"97-93=\", \"17+77=\", \"62-50=\", \"35-0=\", \"12+34=\", \"59-49=\", \"58-43=\", \"57+1=\", \"58-51=\", \"67+17=\", \"20+51=\", \"97-24=\", \"28+18=\", \"77+15=\", \"84-65=\", \"69-34=\", \"29+65=\", \"99-54=\", \"41-35=\", \"87-54=\", \"66-42=\", \"37+1=\", \"23+19=\", \"37+13=\", \"9+80=\", \"47+20=\", \"28-9=\", \"31+41=\", \"44+13=\", \"29+53=\", \"16+43=\", \"46-23=\", \"42-16=\", \"43+37=\", \"81-33=\", \"59-23=\", \"42-27=\", \"52+41=\", \"11+83=\", \"72-5=\", \"92-13=\", \"80-37=\", \"89-56=\", \"23+53=\", \"75-13=\", \"52-27=\", \"9+73=\", \"27+43=\", \"71-39=\", \"12+57=\", \"61-60=\", \"5+45=\", \"47-34=\", \"89-83=\", \"91-55=\", \"57+39=\", \"27+26=\", \"4+59=\", \"94-61=\", \"42+45=\", \"80-45=\", \"51+41=\", \"11+72=\", \"75-71=\", \"93+2=\"];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== newValues.length) {\n  throw new Error(\n    `Expected ${newValues.length} paragraphs, found ${items.length}`\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  items[i].insertText(newValues[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# This script updates the heading date line and the 100 arithmetic\n# practice cells laid out in the single 20x5 table that follows it.\n#\n# $newValues[0] is the new text for the heading paragraph (the date\n# line). $newValues[1..100] are the new texts for the table cells, in\n# row-major order (row 1 col 1..5, row 2 col 1..5, ...), matching the\n# single table in the document (20 rows x 5 columns).\n#\n# We address every paragraph/cell strictly by its position (paragraph\n# index / row+column), never by searching for the old text, because a\n# few of the new values collide with other cells' old values (e.g.\n# \"72-5=\" is the old text of one cell and becomes the new text of a\n# different cell), which would make blind find/replace ambiguous.\n\n$newValues = @(\n    '2025-01-02 Thursday',\n    '5+12=',\n    '58-0=',\n    '62+20=',\n    '78-75=',\n    '88+5=',\n    '17+74=',\n    '97-29=',\n    '65-28=',\n    '27+44=',\n    '39+57=',\n    '50-5=',\n    '28+39=',\n    '16+25=',\n    '3+91=',\n    '31-31=',\n    '37+45=',\n    '10+32=',\n    '1+23=',\n    '65-49=',\n    '60-48=',\n    '18+17=',\n    '42+17=',\n    '42+17=',\n    '91+3=',\n    '27+33=',\n    '11+25=',\n    '92-50=',\n    '52+17=',\n    '43+15=',\n    '9+24=',\n    '52-10=',\n    '48-22=',\n    '62+15=',\n    '5+35=',\n    '16+70=',\n    '97-93=',\n    '17+77=',\n    '62-50=',\n    '35-0=',\n    '12+34=',\n    '59-49=',\n    '58-43=',\n    '57+1=',\n    '58-51=',\n    '67+17=',\n    '20+51=',\n    '97-24=',\n    '28+18=',\n    '77+15=',\n    '84-65=',\n    '69-34=',\n    '29+65=',\n    '99-54=',\n    '41-35=',\n    '87-54=',\n    '66-42=',\n    '37+1=',\n    '23+19=',\n    '37+13=',\n    '9+80=',\n    '47+20=',\n    '28-9=',\n    '31+41=',\n    '44+13=',\n    '29+53=',\n    '16+43=',\n    '46-23=',\n    '42-16=',\n    '43+37=',\n    '81-33=',\n    '59-23=',\n    '42-27=',\n    '52+41=',\n    '11+83=',\n    '72-5=',\n    '92-13=',\n    '80-37=',\n    '89-56=',\n    '23+53=',\n    '75-13=',\n    '52-27=',\n    '9+73=',\n    '27+43=',\n    '71-39=',\n    '12+57=',\n    '61-60=',\n    '5+45=',\n    '47-34=',\n    '89-83=',\n    '91-55=',\n    '57+39=',\n    '27+26=',\n    '4+59=',\n    '94-61=',\n    '42+45=',\n    '80-45=',\n    '51+41=',\n    '11+72=',\n    '75-71=',\n    '93+2='\n)\n\n$d = $word.ActiveDocument\n\n# 1) Heading date paragraph (the first paragraph in the document, above\n#    the table).\n$d.Paragraphs(1).Range.Text = $newValues[0]\n\n# 2) Table cells, row by row, column by column.\n$t = $d.Tables(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\nif (($rows * $cols) -ne ($newValues.Count - 1)) {\n  throw \"Expected $($newValues.Count - 1) table cells, found $($rows * $cols)\"\n}\n\n$k = 1\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $t.Cell($r, $c).Range.Text = $newValues[$k]\n    $k = $k + 1\n  }\n}\n"}
